$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "kitchens_1"
$ws.Range("B1").Value = "bedrooms_1"
$ws.Range("C1").Value = "living_rooms_1"
